# Auto-generated edit script applying the numeric corrections described in the commit diff.
# Each block updates the H-N "price/profit" columns for a specific Leve row on a specific sheet,
# matching the authoritative before/after values from the canonical OOXML diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 66724.45
$ws.Range("I113").Value = 3500
$ws.Range("J113").Value = 142593.8
$ws.Range("K113").Value = 3500
$ws.Range("L113").Value = 142593.8
$ws.Range("M113").Value = -246
$ws.Range("N113").Value = -149101.8

$ws.Range("H116").Value = 15480034
$ws.Range("J116").Value = 3186
$ws.Range("L116").Value = 3186
$ws.Range("N116").Value = -10070

$ws.Range("H125").Value = 2383.077
$ws.Range("I125").Value = 1775.5555
$ws.Range("J125").Value = 3750
$ws.Range("K125").Value = 15979.9995
$ws.Range("L125").Value = 33750
$ws.Range("M125").Value = -13519.9995
$ws.Range("N125").Value = -38670

$ws.Range("H132").Value = 4458.892
$ws.Range("I132").Value = 1571.5358
$ws.Range("J132").Value = 13441.777
$ws.Range("K132").Value = 4714.607400000001
$ws.Range("L132").Value = 40325.331
$ws.Range("M132").Value = -2184.607400000001
$ws.Range("N132").Value = -45385.331

$ws.Range("H138").Value = 5917.469
$ws.Range("J138").Value = 6221.311
$ws.Range("L138").Value = 18663.933
$ws.Range("N138").Value = -28943.933

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1091584.1
$ws.Range("I2").Value = 1164309.8
$ws.Range("J2").Value = 700
$ws.Range("K2").Value = 1164309.8
$ws.Range("L2").Value = 700
$ws.Range("M2").Value = -1164196.8
$ws.Range("N2").Value = -926

$ws.Range("H32").Value = 26202.742
$ws.Range("I32").Value = 26142.834
$ws.Range("K32").Value = 26142.834
$ws.Range("M32").Value = -25855.834

$ws.Range("H45").Value = 2478.7058
$ws.Range("I45").Value = 1778.5454
$ws.Range("K45").Value = 1778.5454
$ws.Range("M45").Value = -1401.5454

$ws.Range("H61").Value = 3365.9092
$ws.Range("I61").Value = 2506.0322
$ws.Range("J61").Value = 5416.385
$ws.Range("K61").Value = 2506.0322
$ws.Range("L61").Value = 5416.385
$ws.Range("M61").Value = -2294.0322
$ws.Range("N61").Value = -5840.385

$ws.Range("H110").Value = 928826.3
$ws.Range("I110").Value = 1021434.2
$ws.Range("K110").Value = 1021434.2
$ws.Range("M110").Value = -1019389.2

$ws.Range("H116").Value = 1091584.1
$ws.Range("I116").Value = 1164309.8
$ws.Range("J116").Value = 700
$ws.Range("K116").Value = 1164309.8
$ws.Range("L116").Value = 700
$ws.Range("M116").Value = -1162015.8
$ws.Range("N116").Value = -5288

$ws.Range("H132").Value = 6453.8335
$ws.Range("J132").Value = 9235.333000000001
$ws.Range("L132").Value = 27705.999
$ws.Range("N132").Value = -32765.999

$ws.Range("H136").Value = 3365.9092
$ws.Range("I136").Value = 2506.0322
$ws.Range("J136").Value = 5416.385
$ws.Range("K136").Value = 7518.096600000001
$ws.Range("L136").Value = 16249.155
$ws.Range("M136").Value = -4968.096600000001
$ws.Range("N136").Value = -21349.155

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1091584.1
$ws.Range("I3").Value = 1164309.8
$ws.Range("J3").Value = 700
$ws.Range("K3").Value = 1164309.8
$ws.Range("L3").Value = 700
$ws.Range("M3").Value = -1164195.8
$ws.Range("N3").Value = -928

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6137.4443
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 6137.4443
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 6137.4443
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -6727.4443

$ws.Range("H34").Value = 6137.4443
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 6137.4443
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 6137.4443
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -6541.4443

$ws.Range("H105").Value = 2067325.5
$ws.Range("I105").Value = 4546456.5
$ws.Range("J105").Value = 1383.1666
$ws.Range("K105").Value = 4546456.5
$ws.Range("L105").Value = 1383.1666
$ws.Range("M105").Value = -4544709.5
$ws.Range("N105").Value = -4877.1666

$ws.Range("H132").Value = 45986060
$ws.Range("I132").Value = 57976270
$ws.Range("K132").Value = 173928810
$ws.Range("M132").Value = -173926280

$ws.Range("H141").Value = 110348.54
$ws.Range("J141").Value = 110841.98
$ws.Range("L141").Value = 110841.98
$ws.Range("N141").Value = -121201.98

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 777.26086
$ws.Range("J5").Value = 1076
$ws.Range("L5").Value = 3228
$ws.Range("N5").Value = -3452

$ws.Range("H8").Value = 205
$ws.Range("I8").Value = 205
$ws.Range("K8").Value = 615
$ws.Range("M8").Value = -476

$ws.Range("H29").Value = 195.71428
$ws.Range("I29").Value = 207.75
$ws.Range("J29").Value = 179.66667
$ws.Range("K29").Value = 623.25
$ws.Range("L29").Value = 539.00001
$ws.Range("M29").Value = -346.25
$ws.Range("N29").Value = -1093.00001

$ws.Range("H32").Value = 600
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()

$ws.Range("H69").Value = 5226.778
$ws.Range("J69").Value = 5255.125
$ws.Range("L69").Value = 15765.375
$ws.Range("N69").Value = -17387.375

$ws.Range("H72").Value = 5226.778
$ws.Range("J72").Value = 5255.125
$ws.Range("L72").Value = 47296.125
$ws.Range("N72").Value = -55408.125

$ws.Range("H87").Value = 4278.5
$ws.Range("I87").Value = 3999.3333
$ws.Range("K87").Value = 11997.9999
$ws.Range("M87").Value = -10749.9999

$ws.Range("H90").Value = 4278.5
$ws.Range("I90").Value = 3999.3333
$ws.Range("K90").Value = 35993.9997
$ws.Range("M90").Value = -29753.9997

$ws.Range("H129").Value = 1344.3793
$ws.Range("J129").Value = 3332.1667
$ws.Range("L129").Value = 9996.500100000001
$ws.Range("N129").Value = -19996.5001

$ws.Range("H135").Value = 777.26086
$ws.Range("J135").Value = 1076
$ws.Range("L135").Value = 9684
$ws.Range("N135").Value = -14754

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 4718081
$ws.Range("I113").Value = 4718081
$ws.Range("K113").Value = 4718081
$ws.Range("M113").Value = -4715911

$ws.Range("H123").Value = 33156.6
$ws.Range("J123").Value = 33156.6
$ws.Range("L123").Value = 33156.6
$ws.Range("N123").Value = -38056.6

$ws.Range("H132").Value = 3356.7778
$ws.Range("I132").Value = 2937.0571
$ws.Range("K132").Value = 8811.1713
$ws.Range("M132").Value = -6281.1713

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 6157.1
$ws.Range("I61").Value = 6734.625
$ws.Range("K61").Value = 6734.625
$ws.Range("M61").Value = -6532.625

$ws.Range("H113").Value = 6157.1
$ws.Range("I113").Value = 6734.625
$ws.Range("K113").Value = 6734.625
$ws.Range("M113").Value = -4564.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 4852.875
$ws.Range("I136").Value = 3411.889
$ws.Range("K136").Value = 10235.667
$ws.Range("M136").Value = -7685.667000000001

$ws.Range("H138").Value = 124185.6
$ws.Range("J138").Value = 132807.25
$ws.Range("L138").Value = 132807.25
$ws.Range("N138").Value = -143087.25

$ws.Range("H140").Value = 79412.75
$ws.Range("J140").Value = 79412.75
$ws.Range("L140").Value = 79412.75
$ws.Range("N140").Value = -89772.75
